$wb = $excel.ActiveWorkbook

$oldMd    = "44fdf768-7dd0-439e-8e50-40cdacd40584.md"
$newMd    = "df1d8dcc-f94d-4e09-8494-ce933d5477c1.md"

$oldDate1 = "2016-03-22 11:02:24"
$newDate1 = "2016-03-22 11:03:02"

$oldZh    = "44fdf768-7dd0-439e-8e50-40cdacd40584.abe27de7eed72fd7d27476ad798ebc7b3747d3b5.zh-cn.xlf"
$newZh    = "df1d8dcc-f94d-4e09-8494-ce933d5477c1.8eea360f748cc614b38cdf996741255eb24a4f5d.zh-cn.xlf"

$oldDate2 = "2016-03-22 11:02:19"
$newDate2 = "2016-03-22 11:02:58"

$oldDe    = "44fdf768-7dd0-439e-8e50-40cdacd40584.abe27de7eed72fd7d27476ad798ebc7b3747d3b5.de-de.xlf"
$newDe    = "df1d8dcc-f94d-4e09-8494-ce933d5477c1.8eea360f748cc614b38cdf996741255eb24a4f5d.de-de.xlf"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newMd
$wsOverview.Range("D2").Value = $newDate1

# zh-cn sheet
$wsZhCn.Range("A2").Value = $newMd
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = $newMd
$wsZhCn.Range("D2").Value = $newZh
$wsZhCn.Hyperlinks.Item(2).TextToDisplay = $newZh
$wsZhCn.Range("E2").Value = $newDate2

# de-de sheet
$wsDeDe.Range("A2").Value = $newMd
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = $newMd
$wsDeDe.Range("D2").Value = $newDe
$wsDeDe.Hyperlinks.Item(2).TextToDisplay = $newDe
$wsDeDe.Range("E2").Value = $newDate1
